# Applies the "Create Client & Corrected Excels" edit:
#  - Summary: B2 786.76 -> 836.76, F2 899.4 -> 849.4
#  - Repayment schedule: K3 987.72 -> 937.72, column L width widened,
#    selection moved to D26
#  - Transactions: F3 786.76 -> 836.76, H3 100 -> 50, selection moved to
#    D15, and it is no longer the active/selected tab
#  - Summary becomes the active (selected) sheet/tab instead of Transactions

$wb = $excel.ActiveWorkbook

$wsSummary    = $wb.Worksheets.Item("Summary")
$wsRepayment  = $wb.Worksheets.Item("Repayment schedule")
$wsTransact   = $wb.Worksheets.Item("Transactions")

# --- Summary sheet value updates ---------------------------------------
$wsSummary.Range("B2").Value = 836.76
$wsSummary.Range("F2").Value = 849.4

# --- Repayment schedule sheet updates -----------------------------------
$wsRepayment.Range("K3").Value = 937.72
$wsRepayment.Columns.Item(12).ColumnWidth = 7.5
[void]$wsRepayment.Range("D26").Select()

# --- Transactions sheet updates -----------------------------------------
$wsTransact.Range("F3").Value = 836.76
$wsTransact.Range("H3").Value = 50
[void]$wsTransact.Range("D15").Select()

# --- Active tab / selection: Summary becomes the active sheet ----------
[void]$wsSummary.Activate()
[void]$wsSummary.Range("K13").Select()
